$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numRng = $ws.Range("D2:D51")
$numRng.NumberFormat = "@"

$ws.Range("D2").Value = "43.196.40"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "2.278.44"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "112.75"
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("D6").Value = "265.22"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").Value = "47.41"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "0.0933"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").Value = "9.26"
$ws.Range("E12").Value = "  +7.90%  "
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "15.44"
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").Value = "2.607.95"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "0.864"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "2.302.01"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "43.150.53"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "6.80"
$ws.Range("E20").Value = "  +3.42%  "
$ws.Range("D21").Value = "71.84"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").Value = "2.48"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").Value = "233.85"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "9.58"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "2.87"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D27").Value = "11.37"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").Value = "40.93"
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "173.50"
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("D33").Value = "21.45"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "0.0904"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("D35").Value = "5.70"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "4.65"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("D39").Value = "3.87"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("E40").Value = "  -5.18%  "
$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  +8.31%  "
$ws.Range("D42").Value = "76.54"
$ws.Range("E42").Value = "  +3.41%  "
$ws.Range("D43").Value = "14.13"
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("D45").Value = "6.20"
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "1.38"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("D48").Value = "8.65"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("D49").Value = "103.47"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "1.25"
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("D51").Value = "0.0995"
$ws.Range("E51").Value = "  -0.83%  "

$numRng.ClearFormats()
